$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B..G to C..H)
$ws.Columns("B").Insert()

# New column B should have the same width as column A
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# New header cell for the "Mã khách hàng" column
$ws.Cells.Item(2, 2).Value = "Mã khách hàng"

# New values for the new column
$ws.Cells.Item(4, 2).Value = "D012"
$ws.Cells.Item(5, 2).Value = "KH05"

# Update the active selection
$ws.Range("B6").Select() | Out-Null
